# Slide 1 title: merge "model " + "eXplainability" (flagged err="1") +
# " for retail banking marketing prediction" into a single run reading
# "model explainability for retail banking marketing prediction".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Title
$tr = $shp.TextFrame.TextRange

# Directly assigning the final text would diff it character-by-character
# against the existing runs (which share long common substrings, e.g.
# "plainability" / " for retail banking marketing prediction"), so the
# old run boundaries/formatting (including the err="1" spell-flag) would
# survive. Assigning an unrelated placeholder first breaks that overlap,
# so the following assignment of the real text collapses everything into
# a single fresh run.
$tr.Text = "###"
$tr.Text = "model explainability for retail banking marketing prediction"
